$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Adds two "tarja azul" (blue bar) rectangles to the title slide, matching
# the Office "Colored Fill - Accent 1" shape style (accent1 fill / accent1
# line), centered text, placed above and below the title block.
#
# msoShapeRectangle = 1. Shapes.AddShape(Type, Left, Top, Width, Height)
# takes Left/Top/Width/Height in points (1 pt = 12700 EMU); the values below
# reproduce the target EMU geometry exactly.
# msoThemeColorAccent1 = 5.

$rect1 = $s.Shapes.AddShape(1, 258.2069291338583, 114.20692913385827, 466.75858267716535, 52.13787401574803)
$rect1.Name = "Retângulo 3"
$rect1.Fill.ForeColor.ObjectThemeColor = 5
$rect1.Line.ForeColor.ObjectThemeColor = 5
$rect1.TextFrame.VerticalAnchor = 3
$rect1.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$rect2 = $s.Shapes.AddShape(1, 258.2069291338583, 354.2069291338583, 466.75858267716535, 81.93102362204725)
$rect2.Name = "Retângulo 4"
$rect2.Fill.ForeColor.ObjectThemeColor = 5
$rect2.Line.ForeColor.ObjectThemeColor = 5
$rect2.TextFrame.VerticalAnchor = 3
$rect2.TextFrame.TextRange.ParagraphFormat.Alignment = 2
